# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-19 (row 7 and row 14 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    15 = 2
    16 = 2
    17 = 1
    18 = 2
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
